# Actualización automática desde Streamlit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the dates in C10 and C11
$c10 = $ws.Range("C10").Value2
$c11 = $ws.Range("C11").Value2
$ws.Range("C10").Value = $c11
$ws.Range("C11").Value = $c10

# Remove the last data row (row 38: Consecutivo 37, ABC, 2025-12-15, 789456, FALSE)
$ws.Rows.Item(38).Delete()
